# Build site at 2023-04-12 14:53:07 UTC
# This script fixes the LOM3229.xlsx course-info sheet:
#  - inserts two new rows (for "Docentes responsaveis" values that were
#    missing and had shifted the remaining label/value pairs out of sync)
#  - corrects the content of several B/C cells that, before this fix, were
#    showing the wrong (previous row's) value
#  - fills in the previously-empty Bibliografia cell

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two rows at 13/14 for the professor-name rows -----------------
$ws.Rows("13:14").Insert()

# The row-insert copies column A's bold label formatting down into the new
# A13/A14 cells even though those rows have nothing in column A; clear that
# back out so the cells disappear entirely (matching the other "no-A" rows
# such as row 12/19/24).
$ws.Range("A13:A14").Clear()

# Give the new B/C cells the same look (wrap text, vertical-top, and the
# red "changed value" font in column C) as the rest of the table by
# copying formats from an existing B/C pair before writing the values.
$ws.Range("B3").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C3").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)

$ws.Range("B13").Value = '6495737 - Durval Rodrigues Junior'
$ws.Range("C13").Value = '6495737 - Durval Rodrigues Junior'

$ws.Range("B14").Value = '1643715 - Paulo Atsushi Suzuki'
$ws.Range("C14").Value = '1643715 - Paulo Atsushi Suzuki'

# --- Fix content that had slipped into the wrong row -----------------------

# Objetivos: (row 10) previously showed Durval's name instead of the goal text
$ws.Range("B10").Value = 'Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais.'
$ws.Range("C10").Value = 'Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais.'

# Programa resumido: (row 15, was 13) previously showed the activation date
$ws.Range("B15").Value = 'Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica.'
$ws.Range("C15").Value = 'Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica.'

# Programa: (row 17, was 15) previously showed Durval's name
$ws.Range("B17").Value = 'A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.'
$ws.Range("C17").Value = 'A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.'

# Método: (row 20, was 18) previously showed Suzuki's name
$ws.Range("B20").Value = 'Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.'
$ws.Range("C20").Value = 'Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.'

# Critério: (row 21, was 19) previously showed the "Método" text
$ws.Range("B21").Value = 'Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3'
$ws.Range("C21").Value = 'Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3'

# Norma de recuperação: (row 22, was 20) previously showed the "Critério" text
$ws.Range("B22").Value = 'Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'
$ws.Range("C22").Value = 'Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação'

# Bibliografia: (row 23, was 21) previously showed the "Norma de recuperação" text;
# fill it in with the actual reading list.
$ws.Range("B23").Value = 'PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.
MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.
WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.
CULLITY, B. D.; STOCK, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001.
YACOBI, B. G.; HOLT, D. B.; KAZMERSKI, L. L. Microanalysis of Solids. Plenum Press, New York, 1994.
HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, Wiley, 1999.
HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.'
$ws.Range("C23").Value = 'PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.
MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.
WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.
CULLITY, B. D.; STOCK, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001.
YACOBI, B. G.; HOLT, D. B.; KAZMERSKI, L. L. Microanalysis of Solids. Plenum Press, New York, 1994.
HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, Wiley, 1999.
HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.'
